$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 96.666664
$ws.Range("J55").Value = 90
$ws.Range("L55").Value = 90
$ws.Range("N55").Value = -518
$ws.Range("H62").Value = 3179.6155
$ws.Range("I62").Value = 3666.875
$ws.Range("J62").Value = 2400
$ws.Range("K62").Value = 3666.875
$ws.Range("L62").Value = 2400
$ws.Range("M62").Value = -3042.875
$ws.Range("N62").Value = -3648
$ws.Range("H65").Value = 3179.6155
$ws.Range("I65").Value = 3666.875
$ws.Range("J65").Value = 2400
$ws.Range("K65").Value = 18334.375
$ws.Range("L65").Value = 12000
$ws.Range("M65").Value = -15214.375
$ws.Range("N65").Value = -18240
$ws.Range("H82").Value = 470
$ws.Range("I82").Value = 494.2857
$ws.Range("J82").Value = 300
$ws.Range("K82").Value = 1482.8571
$ws.Range("L82").Value = 900
$ws.Range("M82").Value = -1076.8571
$ws.Range("N82").Value = -1712
$ws.Range("H85").Value = 470
$ws.Range("I85").Value = 494.2857
$ws.Range("J85").Value = 300
$ws.Range("K85").Value = 1482.8571
$ws.Range("L85").Value = 900
$ws.Range("M85").Value = -78.85710000000017
$ws.Range("N85").Value = -3708
$ws.Range("H125").Value = 1220
$ws.Range("I125").Value = 845
$ws.Range("J125").Value = 1520
$ws.Range("K125").Value = 7605
$ws.Range("L125").Value = 13680
$ws.Range("M125").Value = -5145
$ws.Range("N125").Value = -18600
$ws.Range("H132").Value = 2917496.2
$ws.Range("I132").Value = 3324135.2
$ws.Range("J132").Value = 3250.8333
$ws.Range("K132").Value = 9972405.600000001
$ws.Range("L132").Value = 9752.499899999999
$ws.Range("M132").Value = -9969875.600000001
$ws.Range("N132").Value = -14812.4999
$ws.Range("H137").Value = 1237.9429
$ws.Range("I137").Value = 1128.4615
$ws.Range("K137").Value = 3385.3845
$ws.Range("M137").Value = -835.3844999999997

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1642.7142
$ws.Range("I2").Value = 1374.75
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1374.75
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1261.75
$ws.Range("N2").Value = -2226
$ws.Range("H45").Value = 1569.8
$ws.Range("I45").Value = 750
$ws.Range("J45").Value = 1774.75
$ws.Range("K45").Value = 750
$ws.Range("L45").Value = 1774.75
$ws.Range("M45").Value = -373
$ws.Range("N45").Value = -2528.75
$ws.Range("H61").Value = 1848.5862
$ws.Range("I61").Value = 1300.409
$ws.Range("J61").Value = 3571.4285
$ws.Range("K61").Value = 1300.409
$ws.Range("L61").Value = 3571.4285
$ws.Range("M61").Value = -1088.409
$ws.Range("N61").Value = -3995.4285
$ws.Range("H74").Value = 435.5
$ws.Range("I74").Value = 449.7143
$ws.Range("K74").Value = 449.7143
$ws.Range("M74").Value = 424.2857
$ws.Range("H77").Value = 435.5
$ws.Range("I77").Value = 449.7143
$ws.Range("K77").Value = 2248.5715
$ws.Range("M77").Value = 2119.4285
$ws.Range("H116").Value = 1642.7142
$ws.Range("I116").Value = 1374.75
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 1374.75
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 919.25
$ws.Range("N116").Value = -6588
$ws.Range("H132").Value = 4436.436
$ws.Range("I132").Value = 4766.923
$ws.Range("J132").Value = 3775.4614
$ws.Range("K132").Value = 14300.769
$ws.Range("L132").Value = 11326.3842
$ws.Range("M132").Value = -11770.769
$ws.Range("N132").Value = -16386.3842
$ws.Range("H136").Value = 1848.5862
$ws.Range("I136").Value = 1300.409
$ws.Range("J136").Value = 3571.4285
$ws.Range("K136").Value = 3901.227
$ws.Range("L136").Value = 10714.2855
$ws.Range("M136").Value = -1351.227
$ws.Range("N136").Value = -15814.2855

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1642.7142
$ws.Range("I3").Value = 1374.75
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1374.75
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1260.75
$ws.Range("N3").Value = -2228
$ws.Range("H20").Value = 1697.6111
$ws.Range("I20").Value = 1861.7273
$ws.Range("J20").Value = 1439.7142
$ws.Range("K20").Value = 1861.7273
$ws.Range("L20").Value = 1439.7142
$ws.Range("M20").Value = -1614.7273
$ws.Range("N20").Value = -1933.7142
$ws.Range("H80").Value = 346.73685
$ws.Range("J80").Value = 165.86667
$ws.Range("L80").Value = 165.86667
$ws.Range("N80").Value = -2161.86667
$ws.Range("H83").Value = 346.73685
$ws.Range("J83").Value = 165.86667
$ws.Range("L83").Value = 829.33335
$ws.Range("N83").Value = -10813.33335
$ws.Range("H86").Value = 1881.8
$ws.Range("I86").Value = 1423.6
$ws.Range("K86").Value = 1423.6
$ws.Range("M86").Value = -300.5999999999999
$ws.Range("H89").Value = 1881.8
$ws.Range("I89").Value = 1423.6
$ws.Range("K89").Value = 7118
$ws.Range("M89").Value = -1502
$ws.Range("H105").Value = 2043.3572
$ws.Range("I105").Value = 1919.4736
$ws.Range("K105").Value = 1919.4736
$ws.Range("M105").Value = -172.4736
$ws.Range("H107").Value = 1152.875
$ws.Range("J107").Value = 1373.6666
$ws.Range("L107").Value = 1373.6666
$ws.Range("N107").Value = -5213.6666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1542286.4
$ws.Range("I31").Value = 1868.8064
$ws.Range("J31").Value = 3247748.8
$ws.Range("K31").Value = 1868.8064
$ws.Range("L31").Value = 3247748.8
$ws.Range("M31").Value = -1573.8064
$ws.Range("N31").Value = -3248338.8
$ws.Range("H34").Value = 1542286.4
$ws.Range("I34").Value = 1868.8064
$ws.Range("J34").Value = 3247748.8
$ws.Range("K34").Value = 1868.8064
$ws.Range("L34").Value = 3247748.8
$ws.Range("M34").Value = -1666.8064
$ws.Range("N34").Value = -3248152.8
$ws.Range("H132").Value = 1941.7037
$ws.Range("I132").Value = 1156.3334
$ws.Range("J132").Value = 3512.4443
$ws.Range("K132").Value = 3469.0002
$ws.Range("L132").Value = 10537.3329
$ws.Range("M132").Value = -939.0001999999999
$ws.Range("N132").Value = -15597.3329
$ws.Range("H134").Value = 863
$ws.Range("I134").Value = 661.7778
$ws.Range("J134").Value = 1466.6666
$ws.Range("K134").Value = 1985.3334
$ws.Range("L134").Value = 4399.9998
$ws.Range("M134").Value = 549.6666
$ws.Range("N134").Value = -9469.9998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1410
$ws.Range("I64").Value = 890
$ws.Range("K64").Value = 2670
$ws.Range("M64").Value = -2400
$ws.Range("H67").Value = 1410
$ws.Range("I67").Value = 890
$ws.Range("K67").Value = 2670
$ws.Range("M67").Value = -1734
$ws.Range("H86").Value = 500
$ws.Range("I86").Value = 500
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -314
$ws.Range("N86").Value = -3872
$ws.Range("H89").Value = 500
$ws.Range("I89").Value = 500
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 4500
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = 1428
$ws.Range("N89").Value = -16356
$ws.Range("H92").Value = 423
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H121").Value = 1232189.8
$ws.Range("I121").Value = 342
$ws.Range("J121").Value = 1488824.8
$ws.Range("K121").Value = 1026
$ws.Range("L121").Value = 4466474.4
$ws.Range("M121").Value = 284
$ws.Range("N121").Value = -4469094.4
$ws.Range("H129").Value = 1478.6
$ws.Range("I129").Value = 610
$ws.Range("J129").Value = 1695.75
$ws.Range("K129").Value = 1830
$ws.Range("L129").Value = 5087.25
$ws.Range("M129").Value = 3170
$ws.Range("N129").Value = -15087.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6654.4614
$ws.Range("I122").Value = 15100
$ws.Range("J122").Value = 2900.889
$ws.Range("K122").Value = 45300
$ws.Range("L122").Value = 8702.667000000001
$ws.Range("M122").Value = -42850
$ws.Range("N122").Value = -13602.667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 40000
$ws.Range("J87").Value = 40000
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42246
$ws.Range("H90").Value = 40000
$ws.Range("J90").Value = 40000
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -131232
$ws.Range("H100").Value = 3273.5
$ws.Range("I100").Value = 1550
$ws.Range("J100").Value = 4997
$ws.Range("K100").Value = 1550
$ws.Range("L100").Value = 4997
$ws.Range("M100").Value = -1009
$ws.Range("N100").Value = -6079
$ws.Range("H132").Value = 7676.0625
$ws.Range("I132").Value = 8658.6
$ws.Range("J132").Value = 4167
$ws.Range("K132").Value = 25975.8
$ws.Range("L132").Value = 12501
$ws.Range("M132").Value = -23445.8
$ws.Range("N132").Value = -17561
$ws.Range("H133").Value = 43994
$ws.Range("J133").Value = 43994
$ws.Range("L133").Value = 43994
$ws.Range("N133").Value = -49054
$ws.Range("H136").Value = 8856.4375
$ws.Range("I136").Value = 26876
$ws.Range("K136").Value = 80628
$ws.Range("M136").Value = -78078

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H117").Value = 45000
$ws.Range("J117").Value = 45000
$ws.Range("L117").Value = 45000
$ws.Range("N117").Value = -54178
$ws.Range("H136").Value = 3184.4546
$ws.Range("I136").Value = 8681.333000000001
$ws.Range("K136").Value = 26043.999
$ws.Range("M136").Value = -23493.999
